$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data_Type")

# Fix the typo in C4 (Water Movement description)
$ws.Range("C4").Value = "Current velocity, residence times and related measures"

# Update C5 (Water Quality description)
$ws.Range("C5").Value = "Dissolved oxygen, chlorophyll, nutrients, clarity, or multiple metrics"

# Rows 6-14: rows shift up by one (Habitat Analysis row removed), with new IDs/text
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "Temperature"
$ws.Range("C6").Value = "Self explanatory"

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Salinity"
$ws.Range("C7").Value = "Self explanatory"

$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Inundation"
$ws.Range("C8").Value = "Areas that are or will be inundated.  May include depth or other details"

$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Erosion"
$ws.Range("C9").Value = "Rate or risk of shoreline erosion"

$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Sediment"
$ws.Range("C10").Value = "Sediment deposition,  erosion or transport in shallow waters"

$ws.Range("A11").Value = 12
$ws.Range("B11").Value = "Transport"
$ws.Range("C11").Value = "Movement of pollutants, plankton, fish larvae, drifting objects, etc."

$ws.Range("A12").Value = 13
$ws.Range("B12").Value = "Other"
$ws.Range("C12").Value = "Something that does not fit into other categories"

$ws.Range("A13").Value = 14
$ws.Range("B13").Value = "Not Specified"
$ws.Range("C13").Value = ""

$ws.Range("A14").Value = 16
$ws.Range("B14").Value = "Freshwater"
$ws.Range("C14").Value = "Data on freshwater inflows, etc."

# Rows 15-18 (previously Ice, Freshwater, Weather, Biological) are now cleared
$ws.Range("A15:C18").ClearContents()

# Update the named range Data_Type to the new, smaller extent
$wb.Names.Item("Data_Type").RefersTo = "='Data_Type'!`$A`$1:`$C`$14"
